# Apply the PagoMasivoPrestamo update: add two new columns (usuarioAp,
# cuentaActiva) after "archivo", shift the old Estado/Transaccion/Fecha
# columns right, and refresh the sample data row with the new run's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns before the old column E ("Estado") -------------
# This pushes Estado/Transaccion/Fecha from E:G to G:I and (like Excel's
# real Insert behaviour) picks up the formatting of the column to the left,
# so the new cells already end up with the same style as C1:D1 / C2:D2.
$ws.Range("E1:F1").EntireColumn.Insert()

# --- Header row (row 1) ------------------------------------------------
$ws.Range("A1").Value = "usuario"
$ws.Range("B1").Value = "contraseña"
$ws.Range("C1").Value = "tipo_carga"
$ws.Range("D1").Value = "archivo"
$ws.Range("E1").Value = "usuarioAp"
$ws.Range("F1").Value = "cuentaActiva"
$ws.Range("G1").Value = "Estado"
$ws.Range("H1").Value = "Transaccion"
$ws.Range("I1").Value = "Fecha"

# --- Data row (row 2) ----------------------------------------------------
# A leading "'" forces text, matching how this sheet's other data cells
# (B2:D2) are already stored (quote-prefixed / style "3"), and keeps
# numeric-looking values (account/phone numbers) from being reinterpreted
# as numbers.
$ws.Range("A2").Value = "mrobles"
$ws.Range("B2").Value = "'123456"
$ws.Range("C2").Value = "'ARCUPLOAD"
$ws.Range("D2").Value = "'/src/Excel/entregable2/LoanPayment.csv"
$ws.Range("E2").Value = "'walfaro"
$ws.Range("F2").Value = "'1003181948"

# G2:I2 use the plain (unprefixed/unstyled) formatting, same as before.
$ws.Range("G2").Value = "FAILED"
$ws.Range("I2").Value = "17 jul. 2023, 17:12:26"

# H2 ("Transaccion") is blank for this failed run, but still a real
# (empty) text cell rather than a cleared one, so write it via the
# quote-prefix trick and then strip the quote-prefix styling it adds.
$ws.Range("H2").Value = "'"
$ws.Range("H2").Style = "Normal"

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 8.529947916666666   # E -> ~9.36328125
$ws.Columns.Item(6).ColumnWidth = 10.619791666666666  # F -> ~11.453125
$ws.Columns.Item(8).ColumnWidth = 16.893229166666668  # H -> ~17.7265625
$ws.Columns.Item(9).ColumnWidth = 18.166666666666668  # I -> 19.0

# --- Selection --------------------------------------------------------
$ws.Range("H7").Select() | Out-Null
